{"js": "// Update the date line and each two-digit x two-digit multiplication\n// problem in the practice table to its new value. Every source string is\n// unique in the document, so a straightforward ordered find/replace\n// (exact, case-sensitive match) reproduces the edit.\nconst replacements = [\n  [\"2024-02-23 Friday\", \"2024-02-24 Saturday\"],\n  [\"11\u00d773=\", \"84\u00d740=\"],\n  [\"24\u00d742=\", \"25\u00d784=\"],\n  [\"20\u00d720=\", \"96\u00d788=\"],\n  [\"38\u00d781=\", \"67\u00d715=\"],\n  [\"57\u00d719=\", \"68\u00d730=\"],\n  [\"40\u00d727=\", \"56\u00d738=\"],\n  [\"96\u00d722=\", \"72\u00d724=\"],\n  [\"95\u00d712=\", \"80\u00d773=\"],\n  [\"84\u00d745=\", \"80\u00d752=\"],\n  [\"41\u00d762=\", \"14\u00d715=\"],\n  [\"38\u00d735=\", \"99\u00d792=\"],\n  [\"62\u00d761=\", \"12\u00d741=\"],\n  [\"57\u00d788=\", \"74\u00d754=\"],\n  [\"77\u00d775=\", \"65\u00d734=\"],\n  [\"24\u00d790=\", \"39\u00d781=\"],\n  [\"44\u00d748=\", \"55\u00d712=\"],\n  [\"97\u00d778=\", \"42\u00d789=\"],\n  [\"38\u00d724=\", \"62\u00d799=\"],\n  [\"43\u00d759=\", \"25\u00d763=\"],\n  [\"59\u00d785=\", \"93\u00d743=\"],\n  [\"46\u00d722=\", \"43\u00d732=\"],\n  [\"23\u00d799=\", \"84\u00d762=\"],\n  [\"92\u00d758=\", \"47\u00d748=\"],\n  [\"94\u00d797=\", \"38\u00d757=\"],\n  [\"14\u00d735=\", \"29\u00d740=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: \"${oldText}\"`);\n  }\n\n  for (const hit of results.items) {\n    hit.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each two-digit x two-digit multiplication\n# problem in the practice table to the new values. Every source string\n# is unique in the document, so a straightforward ordered\n# find-and-replace (exact, case-sensitive match) reproduces the edit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-02-23 Friday\", \"2024-02-24 Saturday\"),\n  @(\"11\u00d773=\", \"84\u00d740=\"),\n  @(\"24\u00d742=\", \"25\u00d784=\"),\n  @(\"20\u00d720=\", \"96\u00d788=\"),\n  @(\"38\u00d781=\", \"67\u00d715=\"),\n  @(\"57\u00d719=\", \"68\u00d730=\"),\n  @(\"40\u00d727=\", \"56\u00d738=\"),\n  @(\"96\u00d722=\", \"72\u00d724=\"),\n  @(\"95\u00d712=\", \"80\u00d773=\"),\n  @(\"84\u00d745=\", \"80\u00d752=\"),\n  @(\"41\u00d762=\", \"14\u00d715=\"),\n  @(\"38\u00d735=\", \"99\u00d792=\"),\n  @(\"62\u00d761=\", \"12\u00d741=\"),\n  @(\"57\u00d788=\", \"74\u00d754=\"),\n  @(\"77\u00d775=\", \"65\u00d734=\"),\n  @(\"24\u00d790=\", \"39\u00d781=\"),\n  @(\"44\u00d748=\", \"55\u00d712=\"),\n  @(\"97\u00d778=\", \"42\u00d789=\"),\n  @(\"38\u00d724=\", \"62\u00d799=\"),\n  @(\"43\u00d759=\", \"25\u00d763=\"),\n  @(\"59\u00d785=\", \"93\u00d743=\"),\n  @(\"46\u00d722=\", \"43\u00d732=\"),\n  @(\"23\u00d799=\", \"84\u00d762=\"),\n  @(\"92\u00d758=\", \"47\u00d748=\"),\n  @(\"94\u00d797=\", \"38\u00d757=\"),\n  @(\"14\u00d735=\", \"29\u00d740=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $found = $find.Execute(\n    $oldText,  # FindText\n    $true,     # MatchCase\n    $false,    # MatchWholeWord\n    $false,    # MatchWildcards\n    $false,    # MatchSoundsLike\n    $false,    # MatchAllWordForms\n    $true,     # Forward\n    1,         # Wrap (wdFindContinue)\n    $false,    # Format\n    $newText,  # ReplaceWith\n    2          # Replace (wdReplaceAll)\n  )\n\n  if (-not $found) {\n    throw \"Search text not found: $oldText\"\n  }\n}\n"}
